# Auto-generated: apply per-cell text updates (CI data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AF2").Value = "maa://25251 (91.74), ***maa://21730 (26.32), ***maa://39501 (15.15), **maa://36675 (50.0)"
$ws.Range("D3").Value = "maa://40192 (96.55), maa://36987 (96.15), maa://39849 (88.89)"
$ws.Range("H3").Value = "maa://21247 (98.62), *maa://22748 (60.0)"
$ws.Range("L3").Value = "*maa://22880 (64.68), maa://20276 (86.59), *maa://22749 (76.92)"
$ws.Range("T3").Value = "maa://24617 (89.83), **maa://20790 (43.48), ***maa://37170 (16.18), maa://45854 (83.33)"
$ws.Range("D4").Value = "maa://24632 (93.68), **maa://24303 (38.46), maa://22499 (86.67), maa://22746 (100.0)"
$ws.Range("AB4").Value = "*maa://32658 (71.43)"
$ws.Range("D5").Value = "maa://21245 (84.43), maa://22744 (84.62)"
$ws.Range("L5").Value = "*maa://22757 (76.92)"
$ws.Range("L7").Value = "maa://28624 (92.92), maa://24957 (97.78)"
$ws.Range("A8").Value = "更新日期：2025.03.23 13:17:41"
$ws.Range("X9").Value = "maa://26223 (98.0)"
$ws.Range("AB9").Value = "maa://28711 (87.2), ***maa://22740 (5.66), **maa://39938 (46.67), **maa://27377 (42.86), ***maa://25174 (19.05), *maa://45044 (66.67), maa://40166 (96.3)"
$ws.Range("D10").Value = "***maa://25695 (18.52), ***maa://39951 (13.79), ***maa://34206 (19.23), ***maa://39243 (25.0), *maa://45271 (54.76)"
$ws.Range("T10").Value = "maa://27395 (96.57), maa://22755 (87.83), **maa://22756 (40.91), ***maa://21737 (10.61)"
$ws.Range("AB11").Value = "maa://29912 (97.37), maa://22516 (88.37), *maa://20794 (52.24)"
$ws.Range("AB12").Value = "maa://23669 (95.5), maa://36677 (94.03), maa://39872 (92.0)"
$ws.Range("D13").Value = "maa://24999 (92.09), maa://36673 (92.31), maa://25001 (85.92)"
$ws.Range("H13").Value = "*maa://21248 (73.88), **maa://22728 (46.67)"
$ws.Range("O13").Value = "4"
$ws.Range("P13").Value = "maa://22676 (93.02), *maa://22583 (75.36), *maa://22500 (58.7), maa://48321 (100.0)"
$ws.Range("AF13").Value = "**maa://22737 (34.25), maa://39883 (91.89), *maa://39885 (53.33)"
$ws.Range("L14").Value = "maa://26245 (96.77), maa://21288 (96.3), maa://39841 (93.91), maa://36682 (97.44)"
$ws.Range("AB14").Value = "maa://22764 (97.3)"
$ws.Range("D15").Value = "*maa://22743 (78.24), maa://22734 (84.17), *maa://30808 (64.18), *maa://36048 (50.75), maa://45058 (93.33)"
$ws.Range("H15").Value = "maa://24304 (87.95), maa://21478 (89.19)"
$ws.Range("AF15").Value = "maa://21364 (81.12), *maa://36666 (77.59), *maa://22766 (68.07)"
$ws.Range("X16").Value = "maa://28501 (98.1), maa://28051 (96.0)"
$ws.Range("AB16").Value = "maa://26228 (95.15)"
$ws.Range("H17").Value = "maa://22430 (88.83), maa://39599 (84.91)"
$ws.Range("T17").Value = "*maa://42324 (51.52)"
$ws.Range("H18").Value = "maa://24421 (88.67)"
$ws.Range("L18").Value = "maa://22466 (90.64), *maa://22732 (51.58)"
$ws.Range("AF18").Value = "*maa://24313 (59.52), **maa://29784 (50.0), maa://47854 (100.0)"
$ws.Range("AB19").Value = "*maa://30709 (66.29), *maa://36668 (57.5)"
$ws.Range("L20").Value = "maa://41331 (84.94)"
$ws.Range("AF21").Value = "maa://22524 (93.39), *maa://22432 (76.92)"
$ws.Range("X23").Value = "*maa://28503 (68.67)"
$ws.Range("D24").Value = "*maa://24368 (78.5), *maa://46650 (62.5)"
$ws.Range("X24").Value = "maa://29988 (83.91), maa://23504 (93.33), **maa://22892 (40.14), *maa://25141 (77.1), *maa://36663 (77.5), ***maa://22815 (23.08)"
$ws.Range("D25").Value = "maa://29753 (95.24)"
$ws.Range("H25").Value = "*maa://29063 (72.62), *maa://25311 (74.77), ***maa://22725 (4.76), *maa://45047 (62.5)"
$ws.Range("AB26").Value = "maa://42235 (94.64)"
$ws.Range("AF28").Value = "maa://36660 (92.56), *maa://36701 (66.67)"
$ws.Range("H29").Value = "*maa://25175 (66.04)"
$ws.Range("L29").Value = "maa://28432 (93.59), maa://28440 (80.53), maa://31400 (98.82), *maa://28650 (71.43)"
$ws.Range("D30").Value = "maa://45792 (93.75)"
$ws.Range("L31").Value = "maa://35926 (93.38), maa://36258 (84.87), *maa://43904 (72.73)"
$ws.Range("H32").Value = "maa://21895 (97.06), maa://36667 (97.73), **maa://20793 (38.78), maa://22760 (100.0)"
$ws.Range("AF32").Value = "maa://42408 (84.62)"
$ws.Range("T34").Value = "maa://24526 (92.78)"
$ws.Range("AF34").Value = "*maa://32650 (75.0)"
$ws.Range("H39").Value = "maa://36670 (89.22), maa://25199 (84.82), maa://30434 (91.95), maa://45059 (81.82), ***maa://25036 (19.23), *maa://44165 (66.67)"
$ws.Range("P39").Value = "maa://24709 (91.72), maa://47093 (100.0)"
$ws.Range("T39").Value = "maa://45788 (80.77), maa://47079 (92.59), *maa://45790 (73.33)"
$ws.Range("P40").Value = "maa://23278 (95.3), maa://21386 (95.77), maa://36664 (89.29), maa://45550 (100.0)"
$ws.Range("T45").Value = "**maa://39364 (39.47)"
$ws.Range("H47").Value = "maa://27410 (96.5), maa://29661 (97.33), maa://28038 (84.62)"
$ws.Range("H55").Value = "maa://32532 (92.09)"
$ws.Range("H57").Value = "maa://25176 (98.41)"
